$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Match the width Excel gives the new column to the width of column M ("Paid"),
# which is what the newly-inserted column ends up displaying.
$mWidth = $ws.Columns("M:M").ColumnWidth

# Insert a new (blank) column before column N, shifting N/O/P -> O/P/Q.
$ws.Columns("N:N").Insert()
$ws.Columns("N:N").ColumnWidth = $mWidth

# Make "Repayment schedule" the active sheet/tab, with R6 selected.
$ws.Activate() | Out-Null
$ws.Range("R6").Select() | Out-Null
